$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark (end of the "grade_twelve"
#    paragraph). It will be re-created further up the document, right
#    after the date-of-birth line.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Turn the <w:tab/> that sits right after "{year}" (and right
#    before "Noi sinh (Ten Tinh / Thanh pho):") into a plain space,
#    then drop the (collapsed) "_GoBack" bookmark right after it.
# ------------------------------------------------------------------
$yearRange = $d.Content
$yearRange.Find.Execute("{year}")

$tabStart = $yearRange.End
$tabRange = $d.Range($tabStart, $tabStart + 1)
$tabRange.Text = " "

$goBackRange = $d.Range($tabStart + 1, $tabStart + 1)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# ------------------------------------------------------------------
# 3) Rename the "{province}" placeholder to "{place_of_birth2}".
# ------------------------------------------------------------------
$provinceRange = $d.Content
$provinceRange.Find.Execute("{province}", $true, $false, $false, $false, $false, $true, 1, $false, "{place_of_birth2}", 2)
